$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Agency" row (B3) previously held "Oaxaca"; update it to the user's
# email address, matching the existing "Email" row's value and causing
# Excel to de-duplicate the shared string table.
$ws.Range("B3").Value = "zosimo.montiel@sspo.gob.mx"
